$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text number format to D-column cells whose new values look numeric,
# so Excel keeps them as text (matching the original inline-string cell type)
# instead of silently converting "1.00" -> 1, "584.25" -> 584.25 (number), etc.
$textCells = @("D4","D5","D6","D8","D10","D11","D12","D13","D14","D16","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D35","D37","D38","D41","D42","D46","D47","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Update cell values to match the refreshed cryptocurrency snapshot
$ws.Range("D2").Value = "68.316.75"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "3.320.24"
$ws.Range("E3").Value = "  -2.01%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "584.25"
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("D6").Value = "176.32"
$ws.Range("E6").Value = "  -6.48%  "
$ws.Range("D8").Value = "0.582"
$ws.Range("E8").Value = "  -2.73%  "
$ws.Range("D9").Value = "3.314.92"
$ws.Range("E9").Value = "  -1.96%  "
$ws.Range("D10").Value = "0.176"
$ws.Range("E10").Value = "  -4.46%  "
$ws.Range("D11").Value = "0.577"
$ws.Range("E11").Value = "  -2.36%  "
$ws.Range("D12").Value = "45.64"
$ws.Range("E12").Value = "  -4.42%  "
$ws.Range("D13").Value = "0.0000270"
$ws.Range("E13").Value = "  -2.18%  "
$ws.Range("D14").Value = "679.07"
$ws.Range("E14").Value = "  +5.82%  "
$ws.Range("D15").Value = "3.853.60"
$ws.Range("E15").Value = "  -1.96%  "
$ws.Range("D16").Value = "8.37"
$ws.Range("E16").Value = "  -3.25%  "
$ws.Range("D17").Value = "68.398.02"
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").Value = "3.313.15"
$ws.Range("E19").Value = "  -2.19%  "
$ws.Range("D20").Value = "17.46"
$ws.Range("E20").Value = "  -3.61%  "
$ws.Range("D21").Value = "10.92"
$ws.Range("E21").Value = "  -2.59%  "
$ws.Range("D22").Value = "0.888"
$ws.Range("E22").Value = "  -2.84%  "
$ws.Range("D23").Value = "5.47"
$ws.Range("E23").Value = "  +6.69%  "
$ws.Range("D24").Value = "17.19"
$ws.Range("E24").Value = "  -5.38%  "
$ws.Range("D25").Value = "98.66"
$ws.Range("E25").Value = "  -0.98%  "
$ws.Range("D26").Value = "3.88"
$ws.Range("E26").Value = "  -3.86%  "
$ws.Range("D27").Value = "2.68"
$ws.Range("E27").Value = "  -6.75%  "
$ws.Range("D28").Value = "9.26"
$ws.Range("E28").Value = "  -5.11%  "
$ws.Range("D29").Value = "32.97"
$ws.Range("E29").Value = "  +1.32%  "
$ws.Range("D30").Value = "8.42"
$ws.Range("E30").Value = "  -3.72%  "
$ws.Range("D31").Value = "7.04"
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("D32").Value = "592.23"
$ws.Range("E32").Value = "  -2.79%  "
$ws.Range("D33").Value = "10.98"
$ws.Range("E33").Value = "  -1.56%  "
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "3.776.82"
$ws.Range("E34").Value = "  -5.24%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.104"
$ws.Range("E35").Value = "  -2.75%  "
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("D37").Value = "3.42"
$ws.Range("E37").Value = "  -12.56%  "
$ws.Range("D38").Value = "55.57"
$ws.Range("E38").Value = "  -1.08%  "
$ws.Range("E39").Value = "  -0.88%  "
$ws.Range("E40").Value = "  -7.94%  "
$ws.Range("D41").Value = "32.40"
$ws.Range("E41").Value = "  -4.50%  "
$ws.Range("D42").Value = "3.10"
$ws.Range("E42").Value = "  -5.52%  "
$ws.Range("D43").Value = "0.0₃0669"
$ws.Range("E43").Value = "  -5.72%  "
$ws.Range("E44").Value = "  -4.71%  "
$ws.Range("E45").Value = "  -4.56%  "
$ws.Range("D46").Value = "0.0405"
$ws.Range("E46").Value = "  -4.81%  "
$ws.Range("D47").Value = "2.60"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("E48").Value = "  -2.18%  "
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("D50").Value = "1.35"
$ws.Range("E50").Value = "  -2.17%  "
$ws.Range("D51").Value = "2.76"
$ws.Range("E51").Value = "  -2.48%  "
